# Adds the first set of Draugr cards to the Warriors/Mages/Shadow tables.
# The existing "Mages" (row 75) / "Shadow" (row 82) section headers and
# their sub-tables need to shift down by one row to make room for the new
# Warriors entries (rows 63-74), then the new data rows themselves
# (Warriors 63-74, Mages 78-81, Shadow 85-87) are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything from row 75 down by one row so the new Warriors rows
# (63-74) have room before the "Mages" header.
$ws.Rows.Item(75).Insert()

# ---- Warriors table additions (rows 63-74) ----
$warriors = @(
    @(1,  "Rag Draugr",             1,  "light", "yes", "",    "undead"),
    @(2,  "Draugr",                 1,  "light", "yes", "",    "undead"),
    @(3,  "Deathlord",              10, "heavy", "yes", "",    "undead"),
    @(4,  "Scourgelord",            11, "heavy", "yes", "yes", "undead"),
    @(5,  "Draugr Skirmisher",      4,  "heavy", "yes", "",    "undead"),
    @(6,  "Draugr Mob",             3,  "light", "yes", "",    "undead"),
    @(7,  "Fierce Draugr",          5,  "heavy", "yes", "",    "undead"),
    @(8,  "Scourge ",               4,  "heavy", "yes", "",    "undead"),
    @(9,  "Wightlord",              6,  "heavy", "yes", "",    "undead"),
    @(10, "Ancient Nordic Draugr",  5,  "heavy", "yes", "",    "undead"),
    @(11, "Red Eagle",              9,  "heavy", "yes", "",    "undead")
)

$row = 63
foreach ($w in $warriors) {
    $ws.Cells.Item($row, 1).Value = $w[0]
    $ws.Cells.Item($row, 2).Value = $w[1]
    $ws.Cells.Item($row, 3).Value = $w[2]
    $ws.Cells.Item($row, 4).Value = $w[3]
    $ws.Cells.Item($row, 5).Value = $w[4]
    if ($w[5] -ne "") {
        $ws.Cells.Item($row, 6).Value = $w[5]
    }
    $ws.Cells.Item($row, 7).Value = $w[6]
    $row++
}

# Row 74 only has the No./Name columns filled in.
$ws.Cells.Item(74, 1).Value = 12
$ws.Cells.Item(74, 2).Value = "gauldruson brother? One of them is an archer"

# ---- Mages table additions (rows 78-81, after the shift) ----
$mages = @(
    @(1, "Dragon Priest",    13, "frost mage", "yes", "yes", "undead"),
    @(2, "Skeletal Dragon",  12, "frost mage", "yes", "yes", "undead"),
    @(3, "Draugr Overlord",  7,  "frost mage", "yes", "",    "undead"),
    @(4, "Draugr wight",     4,  "spellsword", "yes", "",    "undead")
)

$row = 78
foreach ($m in $mages) {
    $ws.Cells.Item($row, 1).Value = $m[0]
    $ws.Cells.Item($row, 2).Value = $m[1]
    $ws.Cells.Item($row, 3).Value = $m[2]
    $ws.Cells.Item($row, 4).Value = $m[3]
    $ws.Cells.Item($row, 5).Value = $m[4]
    if ($m[5] -ne "") {
        $ws.Cells.Item($row, 6).Value = $m[5]
    }
    $ws.Cells.Item($row, 7).Value = $m[6]
    $row++
}

# ---- Shadow table additions (rows 85-87) ----
$shadow = @(
    @(1, "Draugr Archer",    2, "Archer", "yes"),
    @(2, "Restless Draugr",  3, "Archer", "yes"),
    @(3, "Skeletal Archer",  1, "Archer", "yes")
)

$row = 85
foreach ($s in $shadow) {
    $ws.Cells.Item($row, 1).Value = $s[0]
    $ws.Cells.Item($row, 2).Value = $s[1]
    $ws.Cells.Item($row, 3).Value = $s[2]
    $ws.Cells.Item($row, 4).Value = $s[3]
    $ws.Cells.Item($row, 5).Value = $s[4]
    $row++
}

# Restore the view state to match where the edits were made.
$ws.Range("B75").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 57 | Out-Null
